# Actualización 10 de Mayo
# Updates the statistics on sheets "Estadisticos 1P", "Estadisticos 2P" and
# "Estadisticos Final" after removing a rescatable student record, and
# removes that student's row from the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P (sheet1) ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 11
$ws1.Range("F2").Value = 21
$ws1.Range("G2").Value = 65.63
$ws1.Range("H2").Value = 7.1

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 7
$ws1.Range("F3").Value = 18
$ws1.Range("G3").Value = 72
$ws1.Range("H3").Value = 7.6

# --- Estadisticos 2P (sheet2) ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 13
$ws2.Range("E2").Value = 13
$ws2.Range("F2").Value = 19
$ws2.Range("G2").Value = 59.38
$ws2.Range("H2").Value = 8.300000000000001

$ws2.Range("D3").Value = 7
$ws2.Range("E3").Value = 7
$ws2.Range("F3").Value = 18
$ws2.Range("G3").Value = 72
$ws2.Range("H3").Value = 8.6

# --- Estadisticos Final (sheet3) ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 10
$ws3.Range("F2").Value = 22
$ws3.Range("G2").Value = 68.75
$ws3.Range("H2").Value = 7.3

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 7
$ws3.Range("F3").Value = 18
$ws3.Range("G3").Value = 72
$ws3.Range("H3").Value = 7.6

# --- Rescatables (sheet4): remove the row for ARELLANO NARANJO MARGARITA JAZMIN ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows.Item(2).Delete()
